$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray formatted cell H12 left over from before
$ws.Range("H12").Clear()

# Full target data (rows 2-15), columns: A=index, B=Name, C=Filtername, D=GEN, E=Filtergroup
$data = @(
    @(0,  "Google",               "f_youtube_",        4, "fg_youtube"),
    @(1,  "Google-IPV6",          "f_youtubeipv6_",    4, "fg_youtube"),
    @(2,  "Facebook",             "f_fbspecialip_",    4, "fg_fbspecialip"),
    @(3,  "Facebook-IPV6",        "f_fbspecialipv6_",  4, "fg_fbspecialip"),
    @(4,  "Instagram",            "f_instagram_",      4, "fg_instagram"),
    @(5,  "Instagram-IPV6",       "f_instagramipv6_",  4, "fg_instagram"),
    @(6,  "Tiktok",               "f_tiktokip_",       3, "fg_tiktokip"),
    @(7,  "Tiktok-IPV6",          "f_tiktokipv6_",     3, "fg_tiktokip"),
    @(8,  "Zoom",                 "f_zoomip_",         4, "fg_zoom"),
    @(9,  "Zoom-IPV6",            "f_zoomipv6_",       4, "fg_zoom"),
    @(10, "FreeFacebook",         "f_freefbip_",       4, "fg_freefbip"),
    @(11, "FreeFacebook-IPV6",    "f_freefbipv6_",     4, "fg_freefbip"),
    @(12, "Facebook Header",      "f_fbheaderip_",     4, "fg_fbheaderip"),
    @(13, "Facebook Header-IPV6", "f_fbheaderipv6_",   4, "fg_fbheaderip")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $row++
}

# Column B width changed (now auto best-fit sized wider to fit "Facebook Header-IPV6")
$ws.Columns.Item(2).ColumnWidth = 20.6666666666667

# Update selection to match target
$ws.Range("H22").Select()
